# Auto-generated data refresh script: updates Leve profit-calculation values
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# the latest market-board price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 900.7917
$ws.Range("J28").Value = 974
$ws.Range("L28").Value = 974
$ws.Range("N28").Value = -1944

$ws.Range("H32").Value = 725
$ws.Range("I32").Value = 700
$ws.Range("J32").Value = 750
$ws.Range("K32").Value = 700
$ws.Range("L32").Value = 750
$ws.Range("M32").Value = -374
$ws.Range("N32").Value = -1402

$ws.Range("H70").Value = 12824.25
$ws.Range("I70").Value = 1200
$ws.Range("J70").Value = 13881
$ws.Range("K70").Value = 3600
$ws.Range("L70").Value = 41643
$ws.Range("M70").Value = -3330
$ws.Range("N70").Value = -42183

$ws.Range("H73").Value = 12824.25
$ws.Range("I73").Value = 1200
$ws.Range("J73").Value = 13881
$ws.Range("K73").Value = 3600
$ws.Range("L73").Value = 41643
$ws.Range("M73").Value = -2664
$ws.Range("N73").Value = -43515

$ws.Range("H100").Value = 8092.909
$ws.Range("I100").Value = 3999.75
$ws.Range("J100").Value = 10431.857
$ws.Range("K100").Value = 3999.75
$ws.Range("L100").Value = 10431.857
$ws.Range("M100").Value = -3458.75
$ws.Range("N100").Value = -11513.857

$ws.Range("H106").Value = 911.9231
$ws.Range("I106").Value = 911.9231
$ws.Range("K106").Value = 911.9231
$ws.Range("M106").Value = -280.9231

$ws.Range("H113").Value = 3709.375
$ws.Range("I113").Value = 2668.75
$ws.Range("J113").Value = 4750
$ws.Range("K113").Value = 2668.75
$ws.Range("L113").Value = 4750
$ws.Range("M113").Value = 585.25
$ws.Range("N113").Value = -11258

$ws.Range("H132").Value = 73958.45
$ws.Range("I132").Value = 79354.8
$ws.Range("K132").Value = 238064.4
$ws.Range("M132").Value = -235534.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 22500
$ws.Range("I31").Value = 22500
$ws.Range("K31").Value = 22500
$ws.Range("M31").Value = -22206

$ws.Range("H32").Value = 5436403
$ws.Range("I32").Value = 6667326.5
$ws.Range("J32").Value = 5858
$ws.Range("K32").Value = 6667326.5
$ws.Range("L32").Value = 5858
$ws.Range("M32").Value = -6667039.5
$ws.Range("N32").Value = -6432

$ws.Range("H110").Value = 1357.1666
$ws.Range("I110").Value = 1228.8
$ws.Range("J110").Value = 1999
$ws.Range("K110").Value = 1228.8
$ws.Range("L110").Value = 1999
$ws.Range("M110").Value = 816.2
$ws.Range("N110").Value = -6089

$ws.Range("H122").Value = 1525.9474
$ws.Range("I122").Value = 1250.8125
$ws.Range("J122").Value = 2993.3333
$ws.Range("K122").Value = 3752.4375
$ws.Range("L122").Value = 8979.999899999999
$ws.Range("M122").Value = -1302.4375
$ws.Range("N122").Value = -13879.9999

$ws.Range("H132").Value = 2013.3966
$ws.Range("I132").Value = 1793.7115
$ws.Range("K132").Value = 5381.1345
$ws.Range("M132").Value = -2851.1345

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 922.13043
$ws.Range("J20").Value = 1306.1428
$ws.Range("L20").Value = 1306.1428
$ws.Range("N20").Value = -1800.1428

$ws.Range("H94").Value = 1966.25
$ws.Range("I94").Value = 1622.2222
$ws.Range("J94").Value = 2998.3333
$ws.Range("K94").Value = 1622.2222
$ws.Range("L94").Value = 2998.3333
$ws.Range("M94").Value = -1171.2222
$ws.Range("N94").Value = -3900.3333

$ws.Range("H99").Value = 9235.799999999999
$ws.Range("I99").Value = 5508.579
$ws.Range("J99").Value = 21038.666
$ws.Range("K99").Value = 5508.579
$ws.Range("L99").Value = 21038.666
$ws.Range("M99").Value = -4010.579
$ws.Range("N99").Value = -24034.666

$ws.Range("H105").Value = 1231.2162
$ws.Range("I105").Value = 1274.3462
$ws.Range("K105").Value = 1274.3462
$ws.Range("M105").Value = 472.6538

$ws.Range("H134").Value = 2304.7441
$ws.Range("I134").Value = 1276.2424
$ws.Range("K134").Value = 3828.7272
$ws.Range("M134").Value = -1293.7272

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 81935.34
$ws.Range("I31").Value = 120738.664
$ws.Range("K31").Value = 120738.664
$ws.Range("M31").Value = -120443.664

$ws.Range("H34").Value = 81935.34
$ws.Range("I34").Value = 120738.664
$ws.Range("K34").Value = 120738.664
$ws.Range("M34").Value = -120536.664

$ws.Range("H122").Value = 3742.3125
$ws.Range("I122").Value = 2688.1
$ws.Range("K122").Value = 8064.299999999999
$ws.Range("M122").Value = -5614.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 277841.22
$ws.Range("I2").Value = 416712.66
$ws.Range("J2").Value = 98.333336
$ws.Range("K2").Value = 2500275.96
$ws.Range("L2").Value = 590.000016
$ws.Range("M2").Value = -2500162.96
$ws.Range("N2").Value = -816.000016

$ws.Range("H31").Value = 685.7143
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H36").Value = 2348.875
$ws.Range("I36").Value = 1638.6
$ws.Range("J36").Value = 3532.6667
$ws.Range("K36").Value = 4915.799999999999
$ws.Range("L36").Value = 10598.0001
$ws.Range("M36").Value = -4746.799999999999
$ws.Range("N36").Value = -10936.0001

$ws.Range("H38").Value = 16.5
$ws.Range("I38").Value = 12.5
$ws.Range("J38").Value = 18.5
$ws.Range("K38").Value = 37.5
$ws.Range("L38").Value = 55.5
$ws.Range("M38").Value = 309.5
$ws.Range("N38").Value = -749.5

$ws.Range("H75").Value = 5070.737
$ws.Range("J75").Value = 5656.2666
$ws.Range("L75").Value = 16968.7998
$ws.Range("N75").Value = -18964.7998

$ws.Range("H78").Value = 5070.737
$ws.Range("J78").Value = 5656.2666
$ws.Range("L78").Value = 50906.3994
$ws.Range("N78").Value = -60890.3994

$ws.Range("H87").Value = 16895.889
$ws.Range("I87").Value = 10399.8
$ws.Range("K87").Value = 31199.4
$ws.Range("M87").Value = -29951.4

$ws.Range("H88").Value = 3014
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H90").Value = 16895.889
$ws.Range("I90").Value = 10399.8
$ws.Range("K90").Value = 93598.2
$ws.Range("M90").Value = -87358.2

$ws.Range("H91").Value = 3014
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H111").Value = 3760.8
$ws.Range("I111").Value = 4201
$ws.Range("K111").Value = 12603
$ws.Range("M111").Value = -9536

$ws.Range("H113").Value = 1189.8462
$ws.Range("I113").Value = 1332
$ws.Range("J113").Value = 1147.2
$ws.Range("K113").Value = 3996
$ws.Range("L113").Value = 3441.6
$ws.Range("M113").Value = -1826
$ws.Range("N113").Value = -7781.6

$ws.Range("H119").Value = 3797.8
$ws.Range("I119").Value = 3329.6667
$ws.Range("J119").Value = 4500
$ws.Range("K119").Value = 9989.000100000001
$ws.Range("L119").Value = 13500
$ws.Range("M119").Value = -5151.000100000001
$ws.Range("N119").Value = -23176

$ws.Range("H136").Value = 20000
$ws.Range("I136").Value = 20000
$ws.Range("K136").Value = 60000
$ws.Range("M136").Value = -54900

$ws.Range("H137").Value = 6926.769
$ws.Range("J137").Value = 7087.3335
$ws.Range("L137").Value = 21262.0005
$ws.Range("N137").Value = -31462.0005

$ws.Range("H139").Value = 2858.2
$ws.Range("I139").Value = 1965.6666
$ws.Range("J139").Value = 4197
$ws.Range("K139").Value = 5896.9998
$ws.Range("L139").Value = 12591
$ws.Range("M139").Value = -756.9997999999996
$ws.Range("N139").Value = -22871

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 338.1111
$ws.Range("I2").Value = 219.2
$ws.Range("J2").Value = 486.75
$ws.Range("K2").Value = 219.2
$ws.Range("L2").Value = 486.75
$ws.Range("M2").Value = -106.2
$ws.Range("N2").Value = -712.75

$ws.Range("H102").Value = 3238.1072
$ws.Range("I102").Value = 3333.2173
$ws.Range("J102").Value = 2800.6
$ws.Range("K102").Value = 3333.2173
$ws.Range("L102").Value = 2800.6
$ws.Range("M102").Value = -1711.2173
$ws.Range("N102").Value = -6044.6

$ws.Range("H107").Value = 2198.7
$ws.Range("I107").Value = 2134.8
$ws.Range("J107").Value = 2262.6
$ws.Range("K107").Value = 2134.8
$ws.Range("L107").Value = 2262.6
$ws.Range("M107").Value = -214.8000000000002
$ws.Range("N107").Value = -6102.6

$ws.Range("H113").Value = 2999
$ws.Range("J113").Value = 2999
$ws.Range("L113").Value = 2999
$ws.Range("N113").Value = -7339

$ws.Range("H126").Value = 3243.3333
$ws.Range("I126").Value = 2806.2307
$ws.Range("K126").Value = 8418.6921
$ws.Range("M126").Value = -5948.6921

$ws.Range("H132").Value = 14353.723
$ws.Range("I132").Value = 13597.76
$ws.Range("J132").Value = 16071.818
$ws.Range("K132").Value = 40793.28
$ws.Range("L132").Value = 48215.454
$ws.Range("M132").Value = -38263.28
$ws.Range("N132").Value = -53275.454

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 555.7179599999999
$ws.Range("I16").Value = 480.3871
$ws.Range("J16").Value = 847.625
$ws.Range("K16").Value = 480.3871
$ws.Range("L16").Value = 847.625
$ws.Range("M16").Value = -310.3871
$ws.Range("N16").Value = -1187.625

$ws.Range("H40").Value = 3418.0625
$ws.Range("I40").Value = 2934
$ws.Range("K40").Value = 2934
$ws.Range("M40").Value = -2798

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 166890.73
$ws.Range("I62").Value = 202944.22
$ws.Range("K62").Value = 202944.22
$ws.Range("M62").Value = -202320.22

$ws.Range("H65").Value = 166890.73
$ws.Range("I65").Value = 202944.22
$ws.Range("K65").Value = 1014721.1
$ws.Range("M65").Value = -1011601.1

